$d = $word.ActiveDocument

# 1. Replace the leading wording of the "Look at replacing the ..." bullet
#    with "Replaced " so the sentence reads "Replaced bottom half of the eye
#    LEDs connection wire from the head with the 3 line servo wire to make
#    it more flexible."
$findRange = $d.Content
$found = $findRange.Find.Execute("Look at replacing the ", $true, $false, $false, $false, $false, $true, 1, $false, "Replaced ", 2)
if (-not $found) {
    throw "Could not find target text to replace"
}

# 2. Locate the end of the just-inserted "Replaced " text so we can drop the
#    _GoBack bookmark (marking the most recent edit position) right there -
#    this also removes it from its previous location near "looked to be",
#    since _GoBack is a single, Word-managed bookmark.
$markerRange = $d.Content
$found2 = $markerRange.Find.Execute("Replaced ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find inserted marker text"
}

$bookmarkRange = $d.Range($markerRange.End, $markerRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
